# Auto-generated edit script: applies numeric corrections to currentAveragePrice /
# LevePrice / LeveProfit columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 900
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 950
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 950
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1602
$ws.Range("H58").Value = 3613.6
$ws.Range("J58").Value = 5999.3335
$ws.Range("L58").Value = 17998.0005
$ws.Range("N58").Value = -18298.0005
$ws.Range("H69").Value = 6666.6665
$ws.Range("J69").Value = 7500
$ws.Range("L69").Value = 22500
$ws.Range("N69").Value = -24248
$ws.Range("H70").Value = 42867.824
$ws.Range("I70").Value = 2458.1667
$ws.Range("J70").Value = 64909.453
$ws.Range("K70").Value = 7374.500100000001
$ws.Range("L70").Value = 194728.359
$ws.Range("M70").Value = -7104.500100000001
$ws.Range("N70").Value = -195268.359
$ws.Range("H72").Value = 6666.6665
$ws.Range("J72").Value = 7500
$ws.Range("L72").Value = 67500
$ws.Range("N72").Value = -76236
$ws.Range("H73").Value = 42867.824
$ws.Range("I73").Value = 2458.1667
$ws.Range("J73").Value = 64909.453
$ws.Range("K73").Value = 7374.500100000001
$ws.Range("L73").Value = 194728.359
$ws.Range("M73").Value = -6438.500100000001
$ws.Range("N73").Value = -196600.359
$ws.Range("H80").Value = 25265.834
$ws.Range("I80").Value = 319
$ws.Range("J80").Value = 150000
$ws.Range("K80").Value = 957
$ws.Range("L80").Value = 450000
$ws.Range("M80").Value = 41
$ws.Range("N80").Value = -451996
$ws.Range("H83").Value = 25265.834
$ws.Range("I83").Value = 319
$ws.Range("J83").Value = 150000
$ws.Range("K83").Value = 2871
$ws.Range("L83").Value = 1350000
$ws.Range("M83").Value = 2121
$ws.Range("N83").Value = -1359984
$ws.Range("H92").Value = 993.7143
$ws.Range("I92").Value = 990
$ws.Range("J92").Value = 1007.3333
$ws.Range("K92").Value = 990
$ws.Range("L92").Value = 1007.3333
$ws.Range("M92").Value = 258
$ws.Range("N92").Value = -3503.3333
$ws.Range("H138").Value = 6164.5312
$ws.Range("I138").Value = 6267.0527
$ws.Range("J138").Value = 6014.6924
$ws.Range("K138").Value = 18801.1581
$ws.Range("L138").Value = 18044.0772
$ws.Range("M138").Value = -13661.1581
$ws.Range("N138").Value = -28324.0772

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17820.328
$ws.Range("I32").Value = 8127.1562
$ws.Range("J32").Value = 27513.5
$ws.Range("K32").Value = 8127.1562
$ws.Range("L32").Value = 27513.5
$ws.Range("M32").Value = -7840.1562
$ws.Range("N32").Value = -28087.5
$ws.Range("H45").Value = 2464.3635
$ws.Range("I45").Value = 1652.25
$ws.Range("K45").Value = 1652.25
$ws.Range("M45").Value = -1275.25
$ws.Range("H61").Value = 2300
$ws.Range("I61").Value = 2300
$ws.Range("K61").Value = 2300
$ws.Range("M61").Value = -2088
$ws.Range("H132").Value = 4301.317
$ws.Range("I132").Value = 1655.4642
$ws.Range("K132").Value = 4966.392599999999
$ws.Range("M132").Value = -2436.392599999999
$ws.Range("H136").Value = 2300
$ws.Range("I136").Value = 2300
$ws.Range("K136").Value = 6900
$ws.Range("M136").Value = -4350

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6807
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6807
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6807
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -9053
$ws.Range("H89").Value = 6807
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6807
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 34035
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -45267
$ws.Range("H105").Value = 4507.6665
$ws.Range("I105").Value = 3883.5334
$ws.Range("K105").Value = 3883.5334
$ws.Range("M105").Value = -2136.5334
$ws.Range("H107").Value = 2482.1667
$ws.Range("J107").Value = 4678
$ws.Range("L107").Value = 4678
$ws.Range("N107").Value = -8518
$ws.Range("H134").Value = 2270.9355
$ws.Range("I134").Value = 1664.4
$ws.Range("K134").Value = 4993.200000000001
$ws.Range("M134").Value = -2458.200000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 287.6
$ws.Range("I22").Value = 146
$ws.Range("K22").Value = 146
$ws.Range("M22").Value = 204
$ws.Range("H31").Value = 3730.7
$ws.Range("J31").Value = 4965.067
$ws.Range("L31").Value = 4965.067
$ws.Range("N31").Value = -5555.067
$ws.Range("H34").Value = 3730.7
$ws.Range("J34").Value = 4965.067
$ws.Range("L34").Value = 4965.067
$ws.Range("N34").Value = -5369.067
$ws.Range("H58").Value = 6896.3
$ws.Range("I58").Value = 4106
$ws.Range("J58").Value = 7593.875
$ws.Range("K58").Value = 4106
$ws.Range("L58").Value = 7593.875
$ws.Range("M58").Value = -3903
$ws.Range("N58").Value = -7999.875
$ws.Range("H62").Value = 58118.375
$ws.Range("I62").Value = 9990.333000000001
$ws.Range("J62").Value = 86995.2
$ws.Range("K62").Value = 9990.333000000001
$ws.Range("L62").Value = 86995.2
$ws.Range("M62").Value = -9366.333000000001
$ws.Range("N62").Value = -88243.2
$ws.Range("H65").Value = 58118.375
$ws.Range("I65").Value = 9990.333000000001
$ws.Range("J65").Value = 86995.2
$ws.Range("K65").Value = 49951.665
$ws.Range("L65").Value = 434976
$ws.Range("M65").Value = -46831.665
$ws.Range("N65").Value = -441216
$ws.Range("H132").Value = 976.05884
$ws.Range("I132").Value = 976.05884
$ws.Range("K132").Value = 2928.17652
$ws.Range("M132").Value = -398.17652
$ws.Range("H134").Value = 2741.1428
$ws.Range("I134").Value = 1908.1765
$ws.Range("J134").Value = 6281.25
$ws.Range("K134").Value = 5724.529500000001
$ws.Range("L134").Value = 18843.75
$ws.Range("M134").Value = -3189.529500000001
$ws.Range("N134").Value = -23913.75
$ws.Range("H136").Value = 6896.3
$ws.Range("I136").Value = 4106
$ws.Range("J136").Value = 7593.875
$ws.Range("K136").Value = 12318
$ws.Range("L136").Value = 22781.625
$ws.Range("M136").Value = -9768
$ws.Range("N136").Value = -27881.625

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14999
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 14999
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 44997
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = -45221
$ws.Range("H22").Value = 3637.25
$ws.Range("I22").Value = 799
$ws.Range("K22").Value = 2397
$ws.Range("M22").Value = -2228
$ws.Range("H27").Value = 3637.25
$ws.Range("I27").Value = 799
$ws.Range("K27").Value = 2397
$ws.Range("M27").Value = -2295
$ws.Range("H117").Value = 2142.818
$ws.Range("J117").Value = 4098.5
$ws.Range("L117").Value = 12295.5
$ws.Range("N117").Value = -19179.5
$ws.Range("H138").Value = 7124.75
$ws.Range("I138").Value = 3250
$ws.Range("J138").Value = 10999.5
$ws.Range("K138").Value = 9750
$ws.Range("L138").Value = 32998.5
$ws.Range("M138").Value = -4610
$ws.Range("N138").Value = -43278.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1412.7778
$ws.Range("I97").Value = 1538.5
$ws.Range("K97").Value = 1538.5
$ws.Range("M97").Value = -1042.5
$ws.Range("H122").Value = 1227665.4
$ws.Range("I122").Value = 504999.5
$ws.Range("K122").Value = 1514998.5
$ws.Range("M122").Value = -1512548.5
$ws.Range("H126").Value = 5443.091
$ws.Range("I126").Value = 4969.75
$ws.Range("J126").Value = 5713.5713
$ws.Range("K126").Value = 14909.25
$ws.Range("L126").Value = 17140.7139
$ws.Range("M126").Value = -12439.25
$ws.Range("N126").Value = -22080.7139

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4999.4
$ws.Range("I7").Value = 4999
$ws.Range("K7").Value = 4999
$ws.Range("M7").Value = -4887
$ws.Range("H68").Value = 3936.6667
$ws.Range("J68").Value = 3840
$ws.Range("L68").Value = 3840
$ws.Range("N68").Value = -5338
$ws.Range("H71").Value = 3936.6667
$ws.Range("J71").Value = 3840
$ws.Range("L71").Value = 19200
$ws.Range("N71").Value = -26688
$ws.Range("H82").Value = 2893
$ws.Range("I82").Value = 2733.8572
$ws.Range("J82").Value = 3450
$ws.Range("K82").Value = 2733.8572
$ws.Range("L82").Value = 3450
$ws.Range("M82").Value = -2372.8572
$ws.Range("N82").Value = -4172
$ws.Range("H85").Value = 2893
$ws.Range("I85").Value = 2733.8572
$ws.Range("J85").Value = 3450
$ws.Range("K85").Value = 2733.8572
$ws.Range("L85").Value = 3450
$ws.Range("M85").Value = -1485.8572
$ws.Range("N85").Value = -5946
$ws.Range("H100").Value = 4187.6
$ws.Range("I100").Value = 1896.5714
$ws.Range("J100").Value = 9533.333000000001
$ws.Range("K100").Value = 1896.5714
$ws.Range("L100").Value = 9533.333000000001
$ws.Range("M100").Value = -1355.5714
$ws.Range("N100").Value = -10615.333
$ws.Range("H122").Value = 5248.75
$ws.Range("I122").Value = 5331.6665
$ws.Range("K122").Value = 15994.9995
$ws.Range("M122").Value = -13544.9995
$ws.Range("H126").Value = 4999.4
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527
$ws.Range("H136").Value = 3789.4
$ws.Range("I136").Value = 2983
$ws.Range("K136").Value = 8949
$ws.Range("M136").Value = -6399

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11600
$ws.Range("I122").Value = 4000.25
$ws.Range("K122").Value = 12000.75
$ws.Range("M122").Value = -9550.75
$ws.Range("H126").Value = 74719.21000000001
$ws.Range("I126").Value = 336833
$ws.Range("K126").Value = 1010499
$ws.Range("M126").Value = -1008029
$ws.Range("H132").Value = 2746.5
$ws.Range("I132").Value = 2085.818
$ws.Range("K132").Value = 6257.454000000001
$ws.Range("M132").Value = -3727.454000000001
